$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77, shifting existing rows 77-80 down to 78-81
$ws.Rows.Item(77).Insert()

# Fill in the new row 77 with the new weekly data point
$ws.Cells.Item(77, 1).Value = 11
$ws.Cells.Item(77, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(77, 3).Value = "Bíobío"
$ws.Cells.Item(77, 4).Value = 44747
$ws.Cells.Item(77, 4).NumberFormat = $ws.Cells.Item(78, 4).NumberFormat
$ws.Cells.Item(77, 5).Value = 8
$ws.Cells.Item(77, 6).Value = 100112012
$ws.Cells.Item(77, 7).Value = "Espinaca"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 150
$ws.Cells.Item(77, 11).Value = 8000
$ws.Cells.Item(77, 12).Value = 8500
$ws.Cells.Item(77, 13).Value = 8233
$ws.Cells.Item(77, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(77, 15).Value = "Región Metropolitana"
$ws.Cells.Item(77, 16).Value = 823
$ws.Cells.Item(77, 17).Value = 10
$ws.Cells.Item(77, 18).Value = "Hortaliza"
